$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Contemporary paragraph" -> "Lyrical " + "paragraph" (two runs)
# ------------------------------------------------------------------
# First swap the word itself, keeping everything in one run.
$d.Content.Find.Execute("Contemporary paragraph", $true, $false, $false, $false, $false, $true, 1, $false, "Lyrical paragraph", 2)

# Now split "Lyrical paragraph" into two runs after "Lyrical ", by
# locating the paragraph and forcing a run boundary with a harmless
# formatting toggle (set then clear) on just the "Lyrical " prefix.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    if ($para.Range.Text.Trim() -eq "Lyrical paragraph") {
        $paraStart = $para.Range.Start
        $prefixRange = $d.Range($paraStart, $paraStart + 8)   # "Lyrical "
        $prefixRange.Font.Bold = $true
        $prefixRange.Font.Bold = $false
        break
    }
}

# ------------------------------------------------------------------
# 2. Remove the trailing empty "ListParagraph" bullet that follows
#    "Enrol now form if needed."
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$lastIndex = $paras.Count
$lastPara = $paras.Item($lastIndex)
if ($lastPara.Range.Text.Trim() -eq "") {
    $prevPara = $paras.Item($lastIndex - 1)
    # Remove from just before the previous paragraph's own end-of-
    # paragraph mark through to the end of the (empty) last paragraph,
    # which deletes the empty paragraph together with its own mark
    # while leaving the previous paragraph's own mark/formatting
    # untouched.
    $deleteRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
    $deleteRange.Delete()
}
